# Final minor edits to metadata descriptions to improve clarity
#
# 1. Remove the now-unused "Reference" column (D) entirely.
# 2. Tweak a couple of variable descriptions (Bag / NDay).
# 3. Rewrite the six "Log of ..." descriptions as "Log-transformed ..."
#    (keeping the original author's small typos, e.g. "Log-tranformed").
# 4. Re-fit the remaining description/units columns now that column D is
#    gone and some text got longer/shorter.
# 5. Leave the cursor/selection on B13, matching the author's final state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete column D ("Reference"), which only ever had a header cell ---
$ws.Columns("D:D").Delete() | Out-Null

# --- 2. Update a couple of descriptions in column B ---
$ws.Range("B7").Value = "32 bags total (4 replicates x 8 treatments); note that a total of five bags are excluded from analysis (and from this dataset), as explained in the methods section of the manuscript"
$ws.Range("B9").Value = "Sampling days numbered 1 through 15 (i.e., order of sampling days rather than actual time intervals)"

# --- 3. Rewrite the "Log of ..." variable descriptions ---
$ws.Range("B27").Value = "Log-transformed TotalDensity (added 1 to TotalDensity with a value of zero to allow for the log transformation)"
$ws.Range("B28").Value = "Log-transformed UninfDensity (added 1 to UninfDensity with a value of zero to allow for the log transformation)"
$ws.Range("B29").Value = "Log-tranformed InfDensity (added 1 to InfDensity with a value of zero to allow for the log transformation)"
$ws.Range("B30").Value = "Log-transformed EdChl (there were no zero values for EdChl)"
$ws.Range("B31").Value = "Log-transformed TP (there were no zero values for TP)"
$ws.Range("B32").Value = "Log-transformed TN (there were no zero values for TN)"

# --- 4. Re-fit column widths for the new layout (B grew, C shrank) ---
$ws.Columns("B:B").ColumnWidth = 111.25
$ws.Columns("C:C").ColumnWidth = 4.75

# --- 5. Match the author's final cursor position/selection ---
$ws.Range("B13").Select() | Out-Null
